$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NetMigration (column C) ballpark figures for non-2025 years
$ws.Range("C3").Value = 75000
$ws.Range("C4").Value = 100000
$ws.Range("C5").Value = 150000
$ws.Range("C6").Value = 100000
$ws.Range("C7").Value = 50000

$ws.Range("C9").Value = 20000
$ws.Range("C10").Value = 25000
$ws.Range("C11").Value = 15000
$ws.Range("C12").Value = 10000
$ws.Range("C13").Value = 8000

$ws.Range("C15").Value = 10000
$ws.Range("C16").Value = 15000
$ws.Range("C17").Value = 10000
$ws.Range("C18").Value = 8000
$ws.Range("C19").Value = 5000

$ws.Range("C21").Value = 5000
$ws.Range("C22").Value = 10000
$ws.Range("C23").Value = 15000
$ws.Range("C24").Value = 5000

# Move the active selection to C1
$ws.Range("C1").Select()
